$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for the specified rows to match
# the "repull data, push all data, mean calculation" commit.
$updates = @{
    4  = 0
    5  = -1
    13 = -3
    17 = -5
    21 = 0
    26 = -1
    28 = 1
    33 = -4
    35 = -1
    36 = -3
    37 = 3
    40 = 1
    44 = -2
    45 = -1
    46 = -2
    47 = -3
    50 = -2
    53 = -1
    54 = 0
    57 = -3
    59 = -1
    61 = 1
    62 = 3
    63 = -6
    64 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 6).Value = $updates[$row]
}
